# Apply updated employee absence data values to rows 2-11
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @{ Row = 2;  A = 32050; B = "Marcos Vinicius Teixeira"; C = "Operacoes";              D = "Viagem de negocios"; E = 1; F = 45082; G = 9392.809999999999 }
    @{ Row = 3;  A = 40397; B = "Aurora Correia";            C = "Engenharia";             D = "Problemas pessoais"; E = 6; F = 45083; G = 3595.64 }
    @{ Row = 4;  A = 13278; B = "Marina Vieira";             C = "Financeiro";             D = "Viagem de negocios"; E = 1; F = 45079; G = 3153.63 }
    @{ Row = 5;  A = 8869;  B = "Aylla Santos";              C = "Atendimento ao Cliente"; D = "Doenca";             E = 2; F = 45085; G = 5261.08 }
    @{ Row = 6;  A = 87339; B = "Gabriel Cardoso";           C = "Financeiro";             D = "Problemas pessoais"; E = 8; F = 45101; G = 8297.76 }
    @{ Row = 7;  A = 34196; B = "Yan Fernandes";             C = "Operacoes";              D = "Consulta medica";   E = 5; F = 45082; G = 3114.37 }
    @{ Row = 8;  A = 86874; B = "Antonella Jesus";           C = "P&D";                    D = "Viagem de negocios"; E = 7; F = 45081; G = 3697.52 }
    @{ Row = 9;  A = 65029; B = "Ana Lívia da Paz";          C = "Financeiro";             D = "Doenca";             E = 8; F = 45103; G = 8422.6 }
    @{ Row = 10; A = 78021; B = "João Miguel Ramos";         C = "Operacoes";              D = "Doenca";             E = 5; F = 45083; G = 2022.13 }
    @{ Row = 11; A = 43998; B = "Mirella Moura";             C = "Financeiro";             D = "Doenca";             E = 3; F = 45101; G = 2271.81 }
)

foreach ($rec in $data) {
    $r = $rec.Row
    $ws.Cells.Item($r, 1).Value = $rec.A
    $ws.Cells.Item($r, 2).Value = $rec.B
    $ws.Cells.Item($r, 3).Value = $rec.C
    $ws.Cells.Item($r, 4).Value = $rec.D
    $ws.Cells.Item($r, 5).Value = $rec.E
    $ws.Cells.Item($r, 6).Value = $rec.F
    $ws.Cells.Item($r, 7).Value = $rec.G
}
